$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")

# Simple F-column ("想去人数" / want-to-go count) value bumps
$ws1.Cells.Item(2, 6).Value = 232
$ws1.Cells.Item(3, 6).Value = 421
$ws1.Cells.Item(4, 6).Value = 159
$ws1.Cells.Item(6, 6).Value = 3867
$ws1.Cells.Item(8, 6).Value = 2548
$ws1.Cells.Item(10, 6).Value = 3135
$ws1.Cells.Item(12, 6).Value = 2311
$ws1.Cells.Item(16, 6).Value = 450
$ws1.Cells.Item(17, 6).Value = 9
$ws1.Cells.Item(20, 6).Value = 345
$ws1.Cells.Item(21, 6).Value = 307
$ws1.Cells.Item(22, 6).Value = 388
$ws1.Cells.Item(23, 6).Value = 660
$ws1.Cells.Item(24, 6).Value = 1409
$ws1.Cells.Item(25, 6).Value = 44
$ws1.Cells.Item(26, 6).Value = 7
$ws1.Cells.Item(27, 6).Value = 1302
$ws1.Cells.Item(28, 6).Value = 130
$ws1.Cells.Item(29, 6).Value = 150
$ws1.Cells.Item(30, 6).Value = 17
$ws1.Cells.Item(37, 6).Value = 7
$ws1.Cells.Item(38, 6).Value = 1125
$ws1.Cells.Item(39, 6).Value = 4
$ws1.Cells.Item(40, 6).Value = 473
$ws1.Cells.Item(42, 6).Value = 1311
$ws1.Cells.Item(43, 6).Value = 168
$ws1.Cells.Item(44, 6).Value = 129
$ws1.Cells.Item(45, 6).Value = 102

# Rows 31-36: source data refreshed. The old cancelled listing
# ("梦游园代号鸢ONLY（取消）", row 31) was dropped and every row below it
# shifted up by one, with a brand-new listing ("配音演员 金弦 专场活动")
# appended at the newly-freed row 36. Column A (the numeric index) is
# never touched by this change, so only columns B..I are rewritten here,
# cell by cell, using the final (post-shift) values for each row.
# Column B holds plain "YYYY-MM-DD" text, not real dates, in the source
# file -- force a text number format first so COM does not silently
# convert the assigned string into a date value.
# Row 31
$ws1.Cells.Item(31, 2).NumberFormat = "@"
$ws1.Cells.Item(31, 2).Value = '2024-09-21'
$ws1.Cells.Item(31, 3).Value = '北京·首届Game 同人Only-神秘玩家'
$ws1.Cells.Item(31, 4).Value = '酒仙桥北路2号院798艺术区706后街1号 北京格瑞斯艺术酒店'
$ws1.Cells.Item(31, 5).Value = '2024.09.21 10:00-09.21 16:30'
$ws1.Cells.Item(31, 6).Value = 3
$ws1.Cells.Item(31, 7).Value = 68
$ws1.Cells.Item(31, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91049'
$ws1.Cells.Item(31, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/V1VIjBIV1724211675935.jpeg'

# Row 32
$ws1.Cells.Item(32, 2).NumberFormat = "@"
$ws1.Cells.Item(32, 2).Value = '2024-09-22'
$ws1.Cells.Item(32, 3).Value = '北京·地狱双ip同人ONLY展'
$ws1.Cells.Item(32, 4).Value = '双桥中路50号院 E50艺术园区'
$ws1.Cells.Item(32, 5).Value = '2024.09.22 10:30-09.22 16:00'
$ws1.Cells.Item(32, 6).Value = 51
$ws1.Cells.Item(32, 7).Value = 105
$ws1.Cells.Item(32, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90931'
$ws1.Cells.Item(32, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/c6ObwO4C1724055713128.jpeg'

# Row 33
$ws1.Cells.Item(33, 2).NumberFormat = "@"
$ws1.Cells.Item(33, 2).Value = '2024-10-01'
$ws1.Cells.Item(33, 3).Value = '北京·IDO动漫游戏嘉年华47th'
$ws1.Cells.Item(33, 4).Value = '亦庄荣昌东街6号 北京亦创国际会展中心'
$ws1.Cells.Item(33, 5).Value = '2024.10.01 09:30-10.02 17:00'
$ws1.Cells.Item(33, 6).Value = 4302
$ws1.Cells.Item(33, 7).Value = 5
$ws1.Cells.Item(33, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83826'
$ws1.Cells.Item(33, 9).Value = '//i0.hdslb.com/bfs/openplatform/202405/JL6boAFV1716882961702.jpeg'

# Row 34
$ws1.Cells.Item(34, 2).NumberFormat = "@"
$ws1.Cells.Item(34, 2).Value = '2024-10-01'
$ws1.Cells.Item(34, 3).Value = '北京·第19届IJOY漫展xCGF游戏节'
$ws1.Cells.Item(34, 4).Value = '天辰东路7号 北京国家会议中心'
$ws1.Cells.Item(34, 5).Value = '2024.10.01 09:00-10.02 17:00'
$ws1.Cells.Item(34, 6).Value = 4038
$ws1.Cells.Item(34, 7).Value = 8.8
$ws1.Cells.Item(34, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84127'
$ws1.Cells.Item(34, 9).Value = '//i0.hdslb.com/bfs/openplatform/202405/iR6rV5311717039317028.jpeg'

# Row 35
$ws1.Cells.Item(35, 2).NumberFormat = "@"
$ws1.Cells.Item(35, 2).Value = '2024-10-01'
$ws1.Cells.Item(35, 3).Value = '北京·第五人格同人only同人3.0'
$ws1.Cells.Item(35, 4).Value = '永外高庄138号 北京大红门国际会展中心'
$ws1.Cells.Item(35, 5).Value = '2024.10.01 10:00-10.01 17:00'
$ws1.Cells.Item(35, 6).Value = 76
$ws1.Cells.Item(35, 7).Value = 60
$ws1.Cells.Item(35, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90653'
$ws1.Cells.Item(35, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/UsBZWtUX1723532208881.jpeg'

# Row 36
$ws1.Cells.Item(36, 2).NumberFormat = "@"
$ws1.Cells.Item(36, 2).Value = '2024-10-01'
$ws1.Cells.Item(36, 3).Value = '北京·配音演员 金弦 专场活动'
$ws1.Cells.Item(36, 4).Value = '亦庄荣昌东街6号 北京亦创国际会展中心'
$ws1.Cells.Item(36, 5).Value = '2024.10.01 10:30-10.01 13:30'
$ws1.Cells.Item(36, 6).Value = 13
$ws1.Cells.Item(36, 7).Value = 268
$ws1.Cells.Item(36, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91068'
$ws1.Cells.Item(36, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/vJRCM3vg1724226523747.jpeg'

# ---- Sheet: 本地生活 ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3, 6).Value = 148
$ws3.Cells.Item(4, 6).Value = 2298

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 148
$ws4.Cells.Item(4, 6).Value = 421
$ws4.Cells.Item(7, 6).Value = 159
$ws4.Cells.Item(9, 6).Value = 3867
$ws4.Cells.Item(11, 6).Value = 2548
$ws4.Cells.Item(13, 6).Value = 3135
$ws4.Cells.Item(15, 6).Value = 2311
$ws4.Cells.Item(18, 6).Value = 9
$ws4.Cells.Item(20, 6).Value = 345
$ws4.Cells.Item(21, 6).Value = 307
$ws4.Cells.Item(22, 6).Value = 388
$ws4.Cells.Item(23, 6).Value = 660
$ws4.Cells.Item(24, 6).Value = 1409
$ws4.Cells.Item(25, 6).Value = 44
$ws4.Cells.Item(26, 6).Value = 1302
$ws4.Cells.Item(27, 6).Value = 130
$ws4.Cells.Item(28, 6).Value = 150
$ws4.Cells.Item(30, 6).Value = 51
$ws4.Cells.Item(32, 6).Value = 4302
$ws4.Cells.Item(34, 6).Value = 7
$ws4.Cells.Item(35, 6).Value = 4
$ws4.Cells.Item(38, 6).Value = 473
$ws4.Cells.Item(43, 6).Value = 1311
$ws4.Cells.Item(44, 6).Value = 168
$ws4.Cells.Item(45, 6).Value = 102

